$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.555.77'
$ws.Range('E2').Value = '  +0.24%  '
$ws.Range('D3').Value = '2.598.95'
$ws.Range('E3').Value = '  -0.75%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '522.49'
$ws.Range('E5').Value = '  +0.41%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.05'
$ws.Range('E6').Value = '  +0.95%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.570'
$ws.Range('E8').Value = '  +0.29%  '
$ws.Range('D9').Value = '2.618.65'
$ws.Range('E9').Value = '  -0.41%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.53'
$ws.Range('E10').Value = '  -0.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.102'
$ws.Range('E11').Value = '  -1.31%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.340'
$ws.Range('E12').Value = '  +0.67%  '
$ws.Range('E13').Value = '  -0.09%  '
$ws.Range('D14').Value = '3.057.07'
$ws.Range('E14').Value = '  -0.76%  '
$ws.Range('D15').Value = '58.422.34'
$ws.Range('E15').Value = '  +0.07%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.44'
$ws.Range('E16').Value = '  -2.49%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000134'
$ws.Range('E17').Value = '  -1.42%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '2.598.50'
$ws.Range('E18').Value = '  -1.57%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '340.79'
$ws.Range('E19').Value = '  +1.03%  '
$ws.Range('E20').Value = '  -0.93%  '
$ws.Range('E21').Value = '  -1.33%  '
$ws.Range('E22').Value = '  +2.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '65.53'
$ws.Range('E24').Value = '  +1.66%  '
$ws.Range('E25').Value = '  +0.96%  '
$ws.Range('E26').Value = '  -2.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.998'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').Value = '2.714.79'
$ws.Range('E28').Value = '  -0.70%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.04'
$ws.Range('E29').Value = '  -1.16%  '
$ws.Range('D30').Value = '0.0₃0750'
$ws.Range('E30').Value = '  -5.10%  '
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.20'
$ws.Range('E32').Value = '  -5.15%  '
$ws.Range('E33').Value = '  -0.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.83'
$ws.Range('E34').Value = '  +0.24%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '149.76'
$ws.Range('E35').Value = '  +0.04%  '
$ws.Range('E36').Value = '  -2.03%  '
$ws.Range('E37').Value = '  -3.93%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.874'
$ws.Range('E38').Value = '  -1.12%  '
$ws.Range('E40').Value = '  +2.43%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '36.07'
$ws.Range('E41').Value = '  -0.77%  '
$ws.Range('E42').Value = '  -2.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.997'
$ws.Range('E43').Value = '  -0.20%  '
$ws.Range('E44').Value = '  -0.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '271.46'
$ws.Range('E46').Value = '  -1.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.67'
$ws.Range('E47').Value = '  +0.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.83'
$ws.Range('E48').Value = '  -1.75%  '
$ws.Range('E49').Value = '  -1.75%  '
$ws.Range('D50').Value = '1.974.17'
$ws.Range('E50').Value = '  -2.70%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.70'
$ws.Range('E51').Value = '  +2.40%  '
